$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that immediately follows
#    the "Play Fortune Dragon Queen Exotic Wilds for Free" H1 heading.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Fortune Dragon Queen Exotic
#    Wilds for Free" right before the very last paragraph (the
#    italic "feature image" paragraph).
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$startPos = $lastPara.Range.Start
$insertRange = $d.Range($startPos, $startPos)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fortune Dragon Queen Exotic Wilds for Free</w:t></w:r></w:p><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($xml)

# ------------------------------------------------------------------
# 3) Swap out the stray duplicate empty run left behind in the final
#    paragraph by the merge above (it previously had a single leading
#    empty run; the XML insert added a second one).
# ------------------------------------------------------------------
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalStart = $finalPara.Range.Start
$dupRange = $d.Range($finalStart, $finalStart + 1)
if ($dupRange.Text -eq "") {
    $dupRange.Delete()
}

# ------------------------------------------------------------------
# 4) Replace the old "Based on the review..." image-prompt text with
#    the new meta-description copy, keeping the italic formatting of
#    the run untouched.
# ------------------------------------------------------------------
$oldText = 'Based on the review of "Fortune Dragon Queen Exotic Wilds", here are the details for the feature image: Title: "Fortune Dragon Queen Exotic Wilds" - Cartoon Style Image Image Description: The main character of the feature image is a happy Maya warrior with glasses. The Maya warrior should be standing proudly with a big smile while holding the slot game on a tablet or smartphone. The background should be decorated with Asian symbols and colors to match the theme of the game. The image should be in cartoon style, with vibrant colors and creative designs that capture the attention of online slot game players. Important elements to include: - A happy Maya warrior with glasses holding the slot game on a tablet or smartphone. - Asian symbols and colors in the background. - Cartoon style image with vibrant colors and creative designs that capture the attention of online slot game players. Overall, the feature image should showcase the excitement and fun of playing "Fortune Dragon Queen Exotic Wilds" while highlighting the Asian-inspired theme and unique gameplay features.'
$newText = 'Read our review of Fortune Dragon Queen Exotic Wilds and get free access to play this online slot game for free today.'

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
